# Remove all old terms for MDR introduction.
# The "end_mdr_introduce_time" row (row 6) is deleted entirely, and the
# remaining "start_mdr_introduce_time" parameter (row 5) is renamed to
# "mdr_introduce_time".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# Delete the entire "end_mdr_introduce_time" row; rows below shift up.
$ws.Rows.Item(6).Delete()

# Rename the remaining MDR introduction time parameter.
$ws.Range("A5").Value = "mdr_introduce_time"

# Reflect the resulting active selection (matches the cell that now
# occupies the position where the deleted row used to be).
$ws.Activate()
$ws.Range("B6").Select()
